# resultFormat.xlsx edit:
#  - rename worksheet "Format" -> "Result"
#  - replace the generic "Столбец1..6" table headers with the real
#    result-report column names
#  - column A ("Status") header/body loses its extra fill formatting so it
#    matches column B's style (wrap + vertical-center only)
#  - a few column widths were nudged
#  - selection left on D12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the sheet ---
$ws.Name = "Result"

# --- new header texts (table column names follow the header cell text) ---
$ws.Range("A1").Value = "Статус"
$ws.Range("B1").Value = "Название компании"
$ws.Range("C1").Value = "ИНН"
$ws.Range("D1").Value = "Начало периода"
$ws.Range("E1").Value = "Окончание периода"
$ws.Range("F1").Value = "Ссылка"

# --- column A (header + first data row) drops its custom fill, becomes
#     identical to column B's plain wrap/vertical-center style ---
$ws.Range("A1:A2").WrapText = $true
$ws.Range("A1:A2").VerticalAlignment = -4108
$ws.Range("A1:A2").Interior.Pattern = -4142

# --- column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 15.333333333333334
$ws.Columns.Item(5).ColumnWidth = 17.333333333333336

# --- leave selection on D12, as in the saved file ---
$ws.Range("D12").Select()
